$d = $word.ActiveDocument

$total = $d.Paragraphs.Count
$startPara = $d.Paragraphs($total - 2)
$endPara = $d.Paragraphs($total)

$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()
